# refactored Flip() to call CheckWinner()
#
# The RANK (Enum) reference list in column R (R14:R26) is rewritten so each
# card rank is paired with its numeric value ("Ace = 1", "Two = 2", ...,
# "King = 13") instead of the old mix of bare rank names / bare numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("R14").Value = "Ace = 1"
$ws.Range("R15").Value = "Two = 2"
$ws.Range("R16").Value = "Three = 3"
$ws.Range("R17").Value = "Four = 4"
$ws.Range("R18").Value = "Five = 5"
$ws.Range("R19").Value = "Six = 6"
$ws.Range("R20").Value = "Seven = 7"
$ws.Range("R21").Value = "Eight = 8"
$ws.Range("R22").Value = "Nine = 9"
$ws.Range("R23").Value = "Ten = 10"
$ws.Range("R24").Value = "Jack = 11"
$ws.Range("R25").Value = "Queen = 12"
$ws.Range("R26").Value = "King = 13"

# Match the author's on-screen view after the edit: scrolled down so row 7
# is at the top, with the freshly-edited RANK list (R13:R26) selected.
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1

$ws.Range("R13:R26").Select()
